$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5322491666369729
$ws.Range("C2").Value = 0.148941578193984
$ws.Range("D2").Value = 0.04285421344134477
$ws.Range("E2").Value = 0.4085644444609926
$ws.Range("F2").Value = 1.32657660840006
$ws.Range("I2").Value = 0.9006924609901361
$ws.Range("K2").Value = 0.6678272695219789
$ws.Range("N2").Value = 1.833030973039406
$ws.Range("B3").Value = 0.4838363878310474
$ws.Range("C3").Value = 0.1334357997532436
$ws.Range("D3").Value = 0.04310138727966706
$ws.Range("E3").Value = 0.3563942191803591
$ws.Range("F3").Value = 1.303922196792826
$ws.Range("I3").Value = 0.8946564304235025
$ws.Range("K3").Value = 0.6046018203161623
$ws.Range("N3").Value = 1.844413391541927
$ws.Range("B4").Value = 0.4543624287551609
$ws.Range("C4").Value = 0.1239830996395312
$ws.Range("D4").Value = 0.04325566482211674
$ws.Range("E4").Value = 0.324472434915549
$ws.Range("F4").Value = 1.290810902224877
$ws.Range("I4").Value = 0.8914314546030866
$ws.Range("K4").Value = 0.5660951004765593
$ws.Range("N4").Value = 1.852011909364464
$ws.Range("B5").Value = 0.4424144269603403
$ws.Range("C5").Value = 0.1201477732944625
$ws.Range("D5").Value = 0.04331916815633541
$ws.Range("E5").Value = 0.3114897092474536
$ws.Range("F5").Value = 1.285667909521351
$ws.Range("I5").Value = 0.8902378216161821
$ws.Range("K5").Value = 0.550481461417121
$ws.Range("N5").Value = 1.855261267521549
$ws.Range("B6").Value = 0.4404342561145995
$ws.Range("C6").Value = 0.1195119202896819
$ws.Range("D6").Value = 0.04332975120333593
$ws.Range("E6").Value = 0.3093354223705234
$ws.Range("F6").Value = 1.284825972328491
$ws.Range("I6").Value = 0.8900468896574978
$ws.Range("K6").Value = 0.5478935266658311
$ws.Range("N6").Value = 1.855810044003206
$ws.Range("B7").Value = 0.4542010397815091
$ws.Range("C7").Value = 0.1239313078258988
$ws.Range("D7").Value = 0.04325651868092084
$ws.Range("E7").Value = 0.3242972443614036
$ws.Range("F7").Value = 1.290740733503512
$ws.Range("I7").Value = 0.891414869232996
$ws.Range("K7").Value = 0.5658842140351226
$ws.Range("N7").Value = 1.852055112811811
$ws.Range("B8").Value = 0.5155039376263915
$ws.Range("C8").Value = 0.143580874032665
$ws.Range("D8").Value = 0.0429389190164704
$ws.Range("E8").Value = 0.3905517577016013
$ws.Range("F8").Value = 1.318599092568832
$ws.Range("I8").Value = 0.8985111103451757
$ws.Range("K8").Value = 0.6459614751505001
$ws.Range("N8").Value = 1.836828849818254
$ws.Range("B9").Value = 0.6377411981186185
$ws.Range("C9").Value = 0.1826703475796592
$ws.Range("D9").Value = 0.04233592847687895
$ws.Range("E9").Value = 0.5214642814764403
$ws.Range("F9").Value = 1.379608874740796
$ws.Range("I9").Value = 0.9162657130987881
$ws.Range("K9").Value = 0.8055281737803739
$ws.Range("N9").Value = 1.811825996508219
$ws.Range("B10").Value = 0.728827088558603
$ws.Range("C10").Value = 0.2117565391003211
$ws.Range("D10").Value = 0.04190484261217531
$ws.Range("E10").Value = 0.6184056745954365
$ws.Range("F10").Value = 1.428388380223964
$ws.Range("I10").Value = 0.9316821761194234
$ws.Range("K10").Value = 0.9243815531403641
$ws.Range("N10").Value = 1.796442202010823
$ws.Range("B11").Value = 0.7705524616295065
$ws.Range("C11").Value = 0.2250745489456278
$ws.Range("D11").Value = 0.04171129564394249
$ws.Range("E11").Value = 0.6627068190713317
$ws.Range("F11").Value = 1.451453276089978
$ws.Range("I11").Value = 0.9392177479239763
$ws.Range("K11").Value = 0.9788196776070777
$ws.Range("N11").Value = 1.790097694592916
$ws.Range("B12").Value = 0.7863951258750603
$ws.Range("C12").Value = 0.2301306288500484
$ws.Range("D12").Value = 0.04163837139386661
$ws.Range("E12").Value = 0.6795142584739295
$ws.Range("F12").Value = 1.460314252152472
$ws.Range("I12").Value = 0.9421469536610942
$ws.Range("K12").Value = 0.9994884738221401
$ws.Range("N12").Value = 1.787789673025387
$ws.Range("B13").Value = 0.782981241367338
$ws.Range("C13").Value = 0.2290411342735013
$ws.Range("D13").Value = 0.04165406057031529
$ws.Range("E13").Value = 0.6758930318468259
$ws.Range("F13").Value = 1.458400226884365
$ws.Range("I13").Value = 0.9415127241689731
$ws.Range("K13").Value = 0.9950346511192265
$ws.Range("N13").Value = 1.788282535383573
$ws.Range("B14").Value = 0.771854999204038
$ws.Range("C14").Value = 0.2254902561120957
$ws.Range("D14").Value = 0.04170528874809865
$ws.Range("E14").Value = 0.6640889302758097
$ws.Range("F14").Value = 1.452179726887422
$ws.Range("I14").Value = 0.9394572164113271
$ws.Range("K14").Value = 0.9805190192854809
$ws.Range("N14").Value = 1.789905915284933
$ws.Range("B15").Value = 0.7650453609039687
$ws.Range("C15").Value = 0.2233169244070439
$ws.Range("D15").Value = 0.04173671539700319
$ws.Range("E15").Value = 0.6568627619730165
$ws.Range("F15").Value = 1.448386034050387
$ws.Range("I15").Value = 0.9382080253555358
$ws.Range("K15").Value = 0.9716348686936556
$ws.Range("N15").Value = 1.7909126050574
$ws.Range("B16").Value = 0.7261061122512444
$ws.Range("C16").Value = 0.2108879509359838
$ws.Range("D16").Value = 0.04191754280628146
$ws.Range("E16").Value = 0.6155147692806793
$ws.Range("F16").Value = 1.426898724295611
$ws.Range("I16").Value = 0.9312002668924535
$ws.Range("K16").Value = 0.9208314375737814
$ws.Range("N16").Value = 1.796870027223136
$ws.Range("B17").Value = 0.702292741992153
$ws.Range("C17").Value = 0.2032856243754679
$ws.Range("D17").Value = 0.04202912927031388
$ws.Range("E17").Value = 0.5902026297404319
$ws.Range("F17").Value = 1.413941732251487
$ws.Range("I17").Value = 0.9270354288371081
$ws.Range("K17").Value = 0.8897608618552511
$ws.Range("N17").Value = 1.800692491028158
$ws.Range("B18").Value = 0.6886231829302574
$ws.Range("C18").Value = 0.1989210904332879
$ws.Range("D18").Value = 0.0420935514088896
$ws.Range("E18").Value = 0.5756626194184378
$ws.Range("F18").Value = 1.406571498063172
$ws.Range("I18").Value = 0.9246890778696226
$ws.Range("K18").Value = 0.871924766194212
$ws.Range("N18").Value = 1.802952565514531
$ws.Range("B19").Value = 0.6839995729816337
$ws.Range("C19").Value = 0.1974447175521163
$ws.Range("D19").Value = 0.04211540495426558
$ws.Range("E19").Value = 0.5707427856592346
$ws.Range("F19").Value = 1.404090165742616
$ws.Range("I19").Value = 0.9239030706131359
$ws.Range("K19").Value = 0.865891734044169
$ws.Range("N19").Value = 1.803728334181017
$ws.Range("B20").Value = 0.7048248910092241
$ws.Range("C20").Value = 0.2040940614683961
$ws.Range("D20").Value = 0.04201722580508882
$ws.Range("E20").Value = 0.5928951798732527
$ws.Range("F20").Value = 1.415312502843051
$ws.Range("I20").Value = 0.9274736921074904
$ws.Range("K20").Value = 0.8930647615580085
$ws.Range("N20").Value = 1.800279215160657
$ws.Range("B21").Value = 0.7751218978558541
$ws.Range("C21").Value = 0.2265328835128742
$ws.Range("D21").Value = 0.04169023180484732
$ws.Range("E21").Value = 0.6675552044232091
$ws.Range("F21").Value = 1.454003389171064
$ws.Range("I21").Value = 0.9400589117189071
$ws.Range("K21").Value = 0.9847811324195845
$ws.Range("N21").Value = 1.789426520291698
$ws.Range("B22").Value = 0.8213110909341594
$ws.Range("C22").Value = 0.2412729701814271
$ws.Range("D22").Value = 0.04147866647858489
$ws.Range("E22").Value = 0.7165349768227713
$ws.Range("F22").Value = 1.480029683914438
$ws.Range("I22").Value = 0.9487252564708655
$ws.Range("K22").Value = 1.045039778078944
$ws.Range("N22").Value = 1.782884750697875
$ws.Range("B23").Value = 0.7966363761141224
$ws.Range("C23").Value = 0.2333989132157512
$ws.Range("D23").Value = 0.04159138649835015
$ws.Range("E23").Value = 0.6903757731970188
$ws.Range("F23").Value = 1.466070957124089
$ws.Range("I23").Value = 0.9440593281664604
$ws.Range("K23").Value = 1.012849327666316
$ws.Range("N23").Value = 1.786325620449375
$ws.Range("B24").Value = 0.7036800408627357
$ws.Range("C24").Value = 0.2037285479497939
$ws.Range("D24").Value = 0.04202260652043233
$ws.Range("E24").Value = 0.5916778399363523
$ws.Range("F24").Value = 1.414692531722196
$ws.Range("I24").Value = 0.9272754035829394
$ws.Range("K24").Value = 0.8915709851204952
$ws.Range("N24").Value = 1.800465862557914
$ws.Range("B25").Value = 0.6044516289914554
$ws.Range("C25").Value = 0.172033157138344
$ws.Range("D25").Value = 0.04249694873438425
$ws.Range("E25").Value = 0.4859270393275779
$ws.Range("F25").Value = 1.362414310711188
$ws.Range("I25").Value = 0.9110483733830534
$ws.Range("K25").Value = 0.7620821232947605
$ws.Range("N25").Value = 1.81806746688487

Write-Output "Updated pl_mw values for Case_5_32 (380 kV case)"
